$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "42.751.95"
Set-TextValue "E2" "  +0.80%  "
Set-TextValue "D3" "2.293.15"
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "315.82"
Set-TextValue "E5" "  -0.56%  "
Set-TextValue "D6" "103.89"
Set-TextValue "E6" "  +0.16%  "
Set-TextValue "E7" "  -1.15%  "
Set-TextValue "E8" "  +0.14%  "
Set-TextValue "E9" "  -1.44%  "
Set-TextValue "D10" "39.40"
Set-TextValue "E10" "  -1.09%  "
Set-TextValue "D11" "0.0905"
Set-TextValue "E11" "  -0.79%  "
Set-TextValue "D12" "8.45"
Set-TextValue "E12" "  +1.09%  "
Set-TextValue "E14" "  +3.04%  "
Set-TextValue "D15" "15.29"
Set-TextValue "E15" "  -0.42%  "
Set-TextValue "D16" "2.641.84"
Set-TextValue "E16" "  -0.34%  "
Set-TextValue "D17" "2.298.36"
Set-TextValue "E17" "  -0.10%  "
Set-TextValue "D18" "42.689.64"
Set-TextValue "E18" "  +0.43%  "
Set-TextValue "D19" "14.66"
Set-TextValue "E19" "  +31.28%  "
Set-TextValue "D20" "7.51"
Set-TextValue "E20" "  -0.06%  "
Set-TextValue "D21" "0.0000105"
Set-TextValue "E21" "  -0.38%  "
Set-TextValue "D22" "74.01"
Set-TextValue "E22" "  +1.21%  "
Set-TextValue "E23" "  -0.78%  "
Set-TextValue "D24" "263.47"
Set-TextValue "E24" "  -5.52%  "
Set-TextValue "E25" "  -3.05%  "
Set-TextValue "E26" "  +0.35%  "
Set-TextValue "E27" "  +0.06%  "
Set-TextValue "D28" "2.34"
Set-TextValue "E28" "  -0.33%  "
Set-TextValue "D29" "6.85"
Set-TextValue "E29" "  +16.35%  "
Set-TextValue "D30" "22.31"
Set-TextValue "E30" "  -2.11%  "
Set-TextValue "D31" "37.31"
Set-TextValue "E31" "  +3.66%  "
Set-TextValue "D32" "166.74"
Set-TextValue "E32" "  +0.99%  "
Set-TextValue "D33" "0.0873"
Set-TextValue "E33" "  -0.48%  "
Set-TextValue "D34" "0.130"
Set-TextValue "E34" "  -4.57%  "
Set-TextValue "D35" "2.58"
Set-TextValue "E35" "  -0.53%  "
Set-TextValue "D36" "0.115"
Set-TextValue "E36" "  -2.94%  "
Set-TextValue "D37" "4.57"
Set-TextValue "E37" "  -0.74%  "
Set-TextValue "E38" "  -6.67%  "
Set-TextValue "D39" "3.71"
Set-TextValue "E39" "  -1.56%  "
Set-TextValue "E40" "  -3.26%  "
Set-TextValue "E41" "  +4.34%  "
Set-TextValue "D42" "69.37"
Set-TextValue "E42" "  -0.83%  "
Set-TextValue "E43" "  +0.25%  "
Set-TextValue "E44" "  -0.07%  "
Set-TextValue "D45" "93.23"
Set-TextValue "E45" "  -3.57%  "
Set-TextValue "D46" "12.27"
Set-TextValue "E46" "  +1.22%  "
Set-TextValue "D47" "113.91"
Set-TextValue "E47" "  +1.68%  "
Set-TextValue "D48" "1.730.41"
Set-TextValue "E48" "  +8.13%  "
Set-TextValue "D49" "80.08"
Set-TextValue "E49" "  -1.37%  "
Set-TextValue "D50" "8.78"
Set-TextValue "E50" "  -1.70%  "
Set-TextValue "E51" "  -0.32%  "

Write-Host "Applied all updates"
